$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 6.87565
$ws.Range("C2").Value = 1.2232
$ws.Range("D2").Value = 35.04514999999999
$ws.Range("E2").Value = 0.0359
$ws.Range("F2").Value = 26.92725
$ws.Range("G2").Value = 26.92725
$ws.Range("K2").Value = 67.2974
$ws.Range("L2").Value = 26.927
$ws.Range("M2").Value = 40.3704
$ws.Range("N2").Value = 7.785999999999999
$ws.Range("O2").Value = 32.5844
$ws.Range("B3").Value = 12.049
$ws.Range("C3").Value = 3.54
$ws.Range("D3").Value = 53.524
$ws.Range("F3").Value = 30.926
$ws.Range("G3").Value = 29.193
$ws.Range("H3").Value = 1.733
$ws.Range("I3").Value = 1.733
$ws.Range("K3").Value = 32.437
$ws.Range("L3").Value = 29.193
$ws.Range("M3").Value = 3.244
$ws.Range("N3").Value = 3.244
$ws.Range("B4").Value = 33.529
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 53.46
$ws.Range("F4").Value = 35.893
$ws.Range("G4").Value = 35.8925808219178
$ws.Range("K4").Value = 35.893
$ws.Range("L4").Value = 35.893

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 6.19875
$ws.Range("C2").Value = 1.20455
$ws.Range("D2").Value = 34.63824999999999
$ws.Range("E2").Value = 0.0359
$ws.Range("F2").Value = 26.70205
$ws.Range("G2").Value = 26.70205
$ws.Range("K2").Value = 229.8882
$ws.Range("L2").Value = 26.70200000000001
$ws.Range("M2").Value = 203.1862
$ws.Range("N2").Value = 8.9472
$ws.Range("O2").Value = 194.2392
$ws.Range("B3").Value = 9.163
$ws.Range("C3").Value = 3.818
$ws.Range("D3").Value = 53.394
$ws.Range("F3").Value = 29.79
$ws.Range("G3").Value = 28.304
$ws.Range("H3").Value = 1.486
$ws.Range("I3").Value = 1.486
$ws.Range("K3").Value = 109.947
$ws.Range("L3").Value = 28.30399999999999
$ws.Range("M3").Value = 81.643
$ws.Range("N3").Value = 7.8102
$ws.Range("O3").Value = 73.83279999999999
$ws.Range("B4").Value = 19.899
$ws.Range("C4").Value = 2.786
$ws.Range("D4").Value = 53.876
$ws.Range("F4").Value = 31.612
$ws.Range("G4").Value = 31.61164383561644
$ws.Range("K4").Value = 35.071
$ws.Range("L4").Value = 31.612
$ws.Range("M4").Value = 3.459000000000001
$ws.Range("N4").Value = 3.459000000000001

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 7.723350000000001
$ws.Range("C2").Value = 1.3063
$ws.Range("D2").Value = 32.7205
$ws.Range("F2").Value = 27.15645
$ws.Range("G2").Value = 27.15645
$ws.Range("K2").Value = 77.2728
$ws.Range("L2").Value = 27.156
$ws.Range("M2").Value = 50.1168
$ws.Range("N2").Value = 5.8988
$ws.Range("O2").Value = 44.218
$ws.Range("B3").Value = 12.049
$ws.Range("C3").Value = 3.54
$ws.Range("D3").Value = 53.524
$ws.Range("F3").Value = 31.738
$ws.Range("G3").Value = 29.193
$ws.Range("H3").Value = 2.544
$ws.Range("I3").Value = 2.544
$ws.Range("K3").Value = 29.193
$ws.Range("L3").Value = 29.193
$ws.Range("B4").Value = 33.529
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 53.46
$ws.Range("F4").Value = 35.893
$ws.Range("G4").Value = 35.8925808219178
$ws.Range("K4").Value = 35.893
$ws.Range("L4").Value = 35.893

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 8.298950000000001
$ws.Range("C2").Value = 1.2959
$ws.Range("D2").Value = 32.83895
$ws.Range("E2").Value = 0.0359
$ws.Range("F2").Value = 27.3412
$ws.Range("G2").Value = 27.3412
$ws.Range("K2").Value = 47.029
$ws.Range("L2").Value = 27.341
$ws.Range("M2").Value = 19.688
$ws.Range("N2").Value = 2.7378
$ws.Range("O2").Value = 16.9502
$ws.Range("B3").Value = 12.049
$ws.Range("C3").Value = 3.54
$ws.Range("D3").Value = 53.524
$ws.Range("F3").Value = 31.738
$ws.Range("G3").Value = 29.193
$ws.Range("H3").Value = 2.544
$ws.Range("I3").Value = 2.544
$ws.Range("K3").Value = 29.193
$ws.Range("L3").Value = 29.193
$ws.Range("B4").Value = 33.529
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 53.46
$ws.Range("F4").Value = 35.893
$ws.Range("G4").Value = 35.8925808219178
$ws.Range("K4").Value = 35.893
$ws.Range("L4").Value = 35.893

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 8.0367
$ws.Range("C2").Value = 1.05765
$ws.Range("D2").Value = 33.68215
$ws.Range("E2").Value = 0.0359
$ws.Range("F2").Value = 27.25005000000001
$ws.Range("G2").Value = 27.25005000000001
$ws.Range("K2").Value = 145.9106
$ws.Range("L2").Value = 27.25
$ws.Range("M2").Value = 118.6606
$ws.Range("N2").Value = 3.2118
$ws.Range("O2").Value = 115.4488
$ws.Range("B3").Value = 12.081
$ws.Range("C3").Value = 2.901
$ws.Range("D3").Value = 52.019
$ws.Range("F3").Value = 31.645
$ws.Range("G3").Value = 29.105
$ws.Range("H3").Value = 2.54
$ws.Range("I3").Value = 2.54
$ws.Range("K3").Value = 42.3792
$ws.Range("L3").Value = 29.105
$ws.Range("M3").Value = 13.2742
$ws.Range("N3").Value = 0.2958
$ws.Range("O3").Value = 12.9786
$ws.Range("B4").Value = 33.529
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 52.947
$ws.Range("F4").Value = 35.882
$ws.Range("G4").Value = 35.88203287671233
$ws.Range("K4").Value = 35.926
$ws.Range("L4").Value = 35.882
$ws.Range("M4").Value = 0.044
$ws.Range("N4").Value = 0.044
$ws.Range("O4").Value = 0

Write-Output "Applied SA_Scenarios_FOB update."
